# The underlying dataset rows (2-20) got reshuffled: each row's record
# (Id, Taxonsorteringsordning, Rodlistade, TaxonId, Artnamn, Vetenskapligt
# namn, Auktor, Ost, Nord, Publik kommentar) now belongs to a different
# spreadsheet row than before, matched up by the record's Id (column A).
# Columns K/L/M/N (Alder-Stadium/Kon/Aktivitet/Metod) and AC (Publik
# kommentar) only exist on bird ("Tretaig hackspett", TaxonId 100109)
# records, so their presence must follow the record as it moves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row number -> source row number (source row's old content moves here)
$rowMap = @{
    2  = 9
    3  = 19
    4  = 6
    5  = 7
    6  = 13
    7  = 4
    8  = 10
    9  = 16
    10 = 12
    11 = 8
    12 = 14
    13 = 15
    14 = 2
    15 = 3
    16 = 17
    17 = 11
    18 = 20
    19 = 18
    20 = 5
}

$colA  = 1   # Id
$colB  = 2   # Taxonsorteringsordning
$colD  = 4   # Rodlistade
$colE  = 5   # TaxonId
$colF  = 6   # Artnamn
$colG  = 7   # Vetenskapligt namn
$colH  = 8   # Auktor
$colK  = 11  # Alder-Stadium
$colL  = 12  # Kon
$colM  = 13  # Aktivitet
$colN  = 14  # Metod
$colQ  = 17  # Ost
$colR  = 18  # Nord
$colAC = 29  # Publik kommentar

# 1) Snapshot every source value BEFORE any writes happen (rows overlap
#    as both read-sources and write-destinations in the permutation).
$snapA  = @{}
$snapB  = @{}
$snapD  = @{}
$snapE  = @{}
$snapF  = @{}
$snapG  = @{}
$snapH  = @{}
$snapQ  = @{}
$snapR  = @{}
$snapAC = @{}

foreach ($r in $rowMap.Keys) {
    $snapA[$r]  = $ws.Cells.Item($r, $colA).Value()
    $snapB[$r]  = $ws.Cells.Item($r, $colB).Value()
    $snapD[$r]  = $ws.Cells.Item($r, $colD).Value()
    $snapE[$r]  = $ws.Cells.Item($r, $colE).Value()
    $snapF[$r]  = $ws.Cells.Item($r, $colF).Value()
    $snapG[$r]  = $ws.Cells.Item($r, $colG).Value()
    $snapH[$r]  = $ws.Cells.Item($r, $colH).Value()
    $snapQ[$r]  = $ws.Cells.Item($r, $colQ).Value()
    $snapR[$r]  = $ws.Cells.Item($r, $colR).Value()
    $snapAC[$r] = $ws.Cells.Item($r, $colAC).Value()
}

# 2) Write each destination row's new content from its mapped source row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]

    $ws.Cells.Item($destRow, $colA).Value = $snapA[$srcRow]
    $ws.Cells.Item($destRow, $colB).Value = $snapB[$srcRow]
    $ws.Cells.Item($destRow, $colD).Value = $snapD[$srcRow]
    $ws.Cells.Item($destRow, $colE).Value = $snapE[$srcRow]
    $ws.Cells.Item($destRow, $colF).Value = $snapF[$srcRow]
    $ws.Cells.Item($destRow, $colG).Value = $snapG[$srcRow]
    $ws.Cells.Item($destRow, $colH).Value = $snapH[$srcRow]
    $ws.Cells.Item($destRow, $colQ).Value = $snapQ[$srcRow]
    $ws.Cells.Item($destRow, $colR).Value = $snapR[$srcRow]

    # Bird records (TaxonId 100109) carry Alder-Stadium/Kon/Aktivitet/Metod
    # (always blank text) plus a Publik kommentar; other species carry
    # neither -- so drive both from the TaxonId that just landed here.
    $isBird = ($snapE[$srcRow] -eq 100109)

    if ($isBird) {
        $ws.Cells.Item($destRow, $colK).Value = "'"
        $ws.Cells.Item($destRow, $colK).Style = "Normal"
        $ws.Cells.Item($destRow, $colL).Value = "'"
        $ws.Cells.Item($destRow, $colL).Style = "Normal"
        $ws.Cells.Item($destRow, $colM).Value = "'"
        $ws.Cells.Item($destRow, $colM).Style = "Normal"
        $ws.Cells.Item($destRow, $colN).Value = "'"
        $ws.Cells.Item($destRow, $colN).Style = "Normal"

        $ws.Cells.Item($destRow, $colAC).Value = $snapAC[$srcRow]
    } else {
        $ws.Cells.Item($destRow, $colK).ClearContents()
        $ws.Cells.Item($destRow, $colL).ClearContents()
        $ws.Cells.Item($destRow, $colM).ClearContents()
        $ws.Cells.Item($destRow, $colN).ClearContents()

        $ws.Cells.Item($destRow, $colAC).ClearContents()
    }
}
